$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after the current last row (row 16) to hold the
# previous week's reading that is being displaced by the new one.
$ws.Rows.Item(17).Insert()

# Row 17 becomes a copy of the "old" row 16 (same date/values it had
# before this week's update).
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = "Vega Monumental Concepción"
$ws.Range("C17").Value = "Bíobío"
$ws.Range("D17").Value = 44461
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100107
$ws.Range("H17").Value = "Otros"
$ws.Range("I17").Value = 100107002
$ws.Range("J17").Value = "Chirimoya"
$ws.Range("K17").Value = "Cultivar IV Región"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 29000
$ws.Range("O17").Value = 30000
$ws.Range("P17").Value = 29500
$ws.Range("Q17").Value = '$/bandeja 10 kilos'
$ws.Range("R17").Value = "Provincia de Limarí"
$ws.Range("S17").Value = 2950
$ws.Range("T17").Value = 10

# Give D17 the same date format as the rest of the date column.
$ws.Range("D17").NumberFormat = $ws.Range("D16").NumberFormat

# Row 16 now reflects the newest week's reading.
$ws.Range("D16").Value = 44491
$ws.Range("M16").Value = 150
$ws.Range("N16").Value = 25000
$ws.Range("O16").Value = 26000
$ws.Range("P16").Value = 25467
$ws.Range("S16").Value = 2547
